# Generate Report for Handback
# Refresh the handback-status report: the "70d93753-..." file's handoff/handback
# timestamps (and the Overview rollup of the latest one) are updated to reflect
# a newer Xliff generation/handback cycle.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsZhCn.Cells.Item(2, 8).Value  = "2016-08-12 14:57:38"
$wsZhCn.Cells.Item(2, 11).Value = "2016-08-12 14:58:16"

# de-de: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsDeDe.Cells.Item(2, 8).Value  = "2016-08-12 14:57:45"
$wsDeDe.Cells.Item(2, 11).Value = "2016-08-12 14:58:26"

# Overview: Latest HO Xliff Generate Date (G2) mirrors the newest handoff time
$wsOverview.Cells.Item(2, 7).Value = "2016-08-12 14:57:45"
